$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Update "Förändrad" (column C) date for the top header rows
#    (2-11) that sit above the reshuffled block handled below.
# ---------------------------------------------------------------
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = 46081
}

# ---------------------------------------------------------------
# 2) Rebuild rows 12-71 (columns A-G) with the new row order and
#    the refreshed "Förändrad" date. Row 71 is a brand-new record
#    that did not exist before.
#    Each line below is one target row, pipe-separated columns
#    A,B,C,D,E,F,G. A value prefixed with "S:" is a string, one
#    prefixed with "N:" is a number, and an empty field means the
#    cell should be left blank.
# ---------------------------------------------------------------
$dataText = @"
S:A 15467-2025|N:45747.52465277778|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:1.4
S:A 11221-2024|N:45371.59097222222|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.3
S:A 5869-2023|N:44958|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:3
S:A 9482-2023|N:44981|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:3.7
S:A 25385-2023|N:45089.33711805556|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.7
S:A 43808-2024|N:45571|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:1.5
S:A 28088-2025|N:45817|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:1
S:A 28093-2025|N:45817|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:10.2
S:A 41803-2025|N:45902.63017361111|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:1.3
S:A 12605-2025|N:45733|N:46081|S:SKÅNE LÄN|S:ESLÖV|S:Sveaskog|N:2.4
S:A 11204-2023|N:44992.60141203704|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.7
S:A 25508-2023|N:45089|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.4
S:A 9686-2023|N:44984|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.8
S:A 61020-2024|N:45645.45321759259|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:4.9
S:A 43943-2025|N:45915.36627314815|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:6.1
S:A 47345-2025|N:45930.57266203704|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.7
S:A 47356-2025|N:45930|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.6
S:A 7409-2025|N:45705.36702546296|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:1.5
S:A 36761-2022|N:44805|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:4.8
S:A 42050-2023|N:45177|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:10.3
S:A 25388-2023|N:45089.34219907408|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.5
S:A 35996-2025|N:45863|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.8
S:A 6545-2023|N:44960|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:4.1
S:A 35997-2025|N:45863|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.9
S:A 13437-2021|N:44273|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:1.3
S:A 61035-2024|N:45645.46634259259|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:7.6
S:A 10940-2024|N:45370.45334490741|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.5
S:A 11732-2025|N:45727|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:4.5
S:A 7016-2025|N:45701|N:46081|S:SKÅNE LÄN|S:ESLÖV|S:Sveaskog|N:2.3
S:A 61179-2025|N:46000.58121527778|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.7
S:A 20160-2024|N:45434|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:4.1
S:A 14488-2024|N:45394|N:46081|S:SKÅNE LÄN|S:ESLÖV|S:Sveaskog|N:1.6
S:A 16670-2025|N:45754.43791666667|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:1.5
S:A 34468-2023|N:45139|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:5.5
S:A 10773-2025|N:45722.47409722222|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:1.2
S:A 47870-2025|N:45932|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:6.8
S:A 62496-2023|N:45268|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:2.6
S:A 10111-2026|N:46074.71814814815|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:5.5
S:A 34170-2023|N:45138|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:5.8
S:A 24115-2022|N:44725|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:4.4
S:A 9849-2026|N:46072|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:2.5
S:A 10819-2026|N:46078|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:6.5
S:A 2691-2026|N:46037.66427083333|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.8
S:A 9484-2023|N:44981|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:1.4
S:A 9487-2023|N:44981|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:1.9
S:A 16733-2023|N:45030|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:2.6
S:A 44496-2023|N:45189|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.5
S:A 34466-2023|N:45139|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:4.3
S:A 33865-2024|N:45520|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:1.3
S:A 5812-2022|N:44596|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.5
S:A 38173-2023|N:45161.44363425926|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:2.8
S:A 15357-2023|N:45019|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:1
S:A 1574-2024|N:45306|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:3.6
S:A 49137-2024|N:45594|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:2.4
S:A 52652-2023|N:45225|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:2.9
S:A 32984-2024|N:45517|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:2.9
S:A 18713-2024|N:45426|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:4.1
S:A 23767-2025|N:45793.47238425926|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:3.2
S:A 23773-2025|N:45793.48923611111|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:4.9
S:A 16792-2024|N:45411|N:46081|S:SKÅNE LÄN|S:ESLÖV||N:0.9
"@

$lines = $dataText -split "`r?`n" | Where-Object { $_.Length -gt 0 }
$rowCount = $lines.Count
$colCount = 7

$arr = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $fields = $lines[$i] -split '\|', -1
    for ($j = 0; $j -lt $colCount; $j++) {
        $field = $fields[$j]
        if ($field -eq "") {
            $arr[$i, $j] = $null
        } elseif ($field.StartsWith("S:")) {
            $arr[$i, $j] = $field.Substring(2)
        } elseif ($field.StartsWith("N:")) {
            $arr[$i, $j] = [double]$field.Substring(2)
        }
    }
}

$ws.Range("A12:G71").Value = $arr

# ---------------------------------------------------------------
# 3) Row 71 is new, so it needs the same per-row formatting that
#    every other data row already carries: date formatting on the
#    "Datum"/"Förändrad" cells, zeroed-out species-count columns
#    (H:Q) and a wrap-text "Artnamn" cell (R) left blank.
# ---------------------------------------------------------------
$ws.Range("B71:C71").NumberFormat = "YYYY-MM-DD"
for ($c = 8; $c -le 17; $c++) {
    $ws.Cells.Item(71, $c).Value = 0
}
$ws.Range("R71").WrapText = $true
$ws.Range("R71").Value = ""
